$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D:E data range to Text format so numeric-looking strings
# (e.g. "161.20", "6.61") are stored verbatim instead of being
# auto-coerced to numbers (which would drop formatting like trailing zeros).
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '64.344.46'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').Value = '3.142.45'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '571.97'
$ws.Range('D6').Value = '163.79'
$ws.Range('E6').Value = '  -3.34%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').Value = '0.573'
$ws.Range('E8').Value = '  -5.79%  '
$ws.Range('D9').Value = '3.155.42'
$ws.Range('E9').Value = '  -0.87%  '
$ws.Range('E10').Value = '  -3.31%  '
$ws.Range('D11').Value = '6.61'
$ws.Range('E11').Value = '  -3.00%  '
$ws.Range('D12').Value = '0.386'
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').Value = '3.695.47'
$ws.Range('E13').Value = '  -1.18%  '
$ws.Range('E14').Value = '  -0.72%  '
$ws.Range('D15').Value = '64.378.69'
$ws.Range('E15').Value = '  -0.31%  '
$ws.Range('D16').Value = '25.24'
$ws.Range('E16').Value = '  -0.84%  '
$ws.Range('D17').Value = '3.145.59'
$ws.Range('E17').Value = '  -1.23%  '
$ws.Range('D18').Value = '0.0000154'
$ws.Range('E18').Value = '  -3.20%  '
$ws.Range('D19').Value = '401.59'
$ws.Range('E19').Value = '  -4.22%  '
$ws.Range('D20').Value = '5.26'
$ws.Range('E20').Value = '  -2.05%  '
$ws.Range('D21').Value = '12.54'
$ws.Range('E21').Value = '  -3.02%  '
$ws.Range('D22').Value = '7.09'
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = '68.55'
$ws.Range('E24').Value = '  -2.59%  '
$ws.Range('D25').Value = '0.485'
$ws.Range('E25').Value = '  -0.79%  '
$ws.Range('E26').Value = '  -4.38%  '
$ws.Range('E27').Value = '  -4.53%  '
$ws.Range('D28').Value = '8.81'
$ws.Range('E28').Value = '  -0.24%  '
$ws.Range('D29').Value = '0.995'
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('D31').Value = '1.81'
$ws.Range('E31').Value = '  -1.66%  '
$ws.Range('D32').Value = '21.23'
$ws.Range('E32').Value = '  -2.47%  '
$ws.Range('D33').Value = '161.20'
$ws.Range('E33').Value = '  +2.26%  '
$ws.Range('D34').Value = '6.28'
$ws.Range('E34').Value = '  -0.81%  '
$ws.Range('D35').Value = '4.84'
$ws.Range('E36').Value = '  -2.53%  '
$ws.Range('D37').Value = '1.34'
$ws.Range('E37').Value = '  -2.16%  '
$ws.Range('E38').Value = '  -1.88%  '
$ws.Range('D39').Value = '2.644.98'
$ws.Range('E39').Value = '  -3.14%  '
$ws.Range('D40').Value = '23.77'
$ws.Range('E40').Value = '  -2.54%  '
$ws.Range('D41').Value = '4.07'
$ws.Range('E41').Value = '  -2.90%  '
$ws.Range('D42').Value = '38.45'
$ws.Range('E42').Value = '  -1.90%  '
$ws.Range('D43').Value = '0.691'
$ws.Range('E43').Value = '  -3.18%  '
$ws.Range('D44').Value = '0.0614'
$ws.Range('E44').Value = '  -1.78%  '
$ws.Range('E45').Value = '  -3.79%  '
$ws.Range('D46').Value = '0.0255'
$ws.Range('E46').Value = '  -3.65%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '21.13'
$ws.Range('E47').Value = '  -2.73%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '286.51'
$ws.Range('E48').Value = '  -2.77%  '
$ws.Range('D49').Value = '0.997'
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('D50').Value = '0.0976'
$ws.Range('E50').Value = '  -1.40%  '
$ws.Range('D51').Value = '10.47'
$ws.Range('E51').Value = '  +0.15%  '

# Restore default (General) formatting so cell styles match the original workbook.
$dataRange.ClearFormats()
